$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 363 (pushes the existing
# row 363.. data down to 365..460), mirroring the two new "Ajo"
# price records added at the top of this block in the source diff.
$ws.Rows.Item(363).Insert()
$ws.Rows.Item(363).Insert()

# --- New row 363 -----------------------------------------------------
$ws.Range("A363").Value = 3
$ws.Range("B363").Value = "Femacal de La Calera"
$ws.Range("C363").Value = "Coquimbo"
$ws.Range("D363").Value = 44736
$ws.Range("E363").Value = 5
$ws.Range("F363").Value = 100112003
$ws.Range("G363").Value = "Ajo"
$ws.Range("H363").Value = "Chino"
$ws.Range("I363").Value = "Primera"
$ws.Range("J363").Value = 82
$ws.Range("K363").Value = 17000
$ws.Range("L363").Value = 17500
$ws.Range("M363").Value = 17244
$ws.Range("N363").Value = "$/caja 10 kilos"
$ws.Range("O363").Value = "China"
$ws.Range("P363").Value = 1724
$ws.Range("Q363").Value = 10
$ws.Range("R363").Value = "Hortaliza"

# --- New row 364 -----------------------------------------------------
$ws.Range("A364").Value = 3
$ws.Range("B364").Value = "Femacal de La Calera"
$ws.Range("C364").Value = "Coquimbo"
$ws.Range("D364").Value = 44736
$ws.Range("E364").Value = 5
$ws.Range("F364").Value = 100112003
$ws.Range("G364").Value = "Ajo"
$ws.Range("H364").Value = "Chino"
$ws.Range("I364").Value = "Primera"
$ws.Range("J364").Value = 38
$ws.Range("K364").Value = 21000
$ws.Range("L364").Value = 21000
$ws.Range("M364").Value = 21000
$ws.Range("N364").Value = "$/malla 10 kilos"
$ws.Range("O364").Value = "China"
$ws.Range("P364").Value = 2100
$ws.Range("Q364").Value = 10
$ws.Range("R364").Value = "Hortaliza"
